$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 data (appended log entry for a new automatic test-sync run)
$ws.Range("A6").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$ws.Range("B6").Value = "Beste afzender,`nBedankt voor je e-mail. Helaas kan ik je niet helpen met het plaatsen van bestellingen via e-mail. Je kunt onze webshop bezoeken om de gewenste M5-bouten te bestellen. Mocht je nog vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam] - E-mailassistent"
$ws.Range("C6").Value = "Wil je 100 stuks M5-bouten bestellen?"
$ws.Range("D6").Value = "mailmind.test@zohomail.eu"
$ws.Range("E6").Value = "Bestelling / Levering"
$ws.Range("F6").Value = "2025-07-31 21:29:47"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Ja"
$ws.Range("J6").Value = "Nee"

$wb.Save()
